$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3719
$ws1.Range("F5").Value = 3719
$ws1.Range("F6").Value = 286
$ws1.Range("F7").Value = 5252
$ws1.Range("F8").Value = 572
$ws1.Range("F9").Value = 405
$ws1.Range("F11").Value = 1032
$ws1.Range("F13").Value = 125
$ws1.Range("F16").Value = 348
$ws1.Range("F19").Value = 168
$ws1.Range("F22").Value = 6006
$ws1.Range("F26").Value = 6314
$ws1.Range("F29").Value = 3249
$ws1.Range("F30").Value = 363
$ws1.Range("F31").Value = 740
$ws1.Range("F32").Value = 4453
$ws1.Range("F36").Value = 1106
$ws1.Range("F37").Value = 97
$ws1.Range("F40").Value = 909
$ws1.Range("F41").Value = 1095
$ws1.Range("F42").Value = 2051

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1142

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1142
$ws4.Range("F7").Value = 3719
$ws4.Range("F8").Value = 3719
$ws4.Range("F9").Value = 286
$ws4.Range("F10").Value = 5252
$ws4.Range("F11").Value = 572
$ws4.Range("F12").Value = 405
$ws4.Range("F14").Value = 1032
$ws4.Range("F16").Value = 125
$ws4.Range("F19").Value = 348
$ws4.Range("F23").Value = 168
$ws4.Range("F26").Value = 6006
$ws4.Range("F30").Value = 6314
$ws4.Range("F33").Value = 3249
$ws4.Range("F34").Value = 363
$ws4.Range("F35").Value = 740
$ws4.Range("F36").Value = 4453
$ws4.Range("F41").Value = 1106
$ws4.Range("F42").Value = 97
$ws4.Range("F45").Value = 909
$ws4.Range("F46").Value = 1095
$ws4.Range("F48").Value = 2051
